# Update cryptos list (Price/Volume(1h) columns) with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.789.37"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "2.905.71"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("D5").Value = "'526.79"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").Value = "'143.73"
$ws.Range("E6").Value = "  -5.59%  "
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").Value = "2.914.96"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("E10").Value = "  -4.58%  "
$ws.Range("D11").Value = "'6.05"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").Value = "3.411.91"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "60.764.43"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "'22.52"
$ws.Range("E16").Value = "  -5.93%  "
$ws.Range("D17").Value = "2.914.37"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("D21").Value = "'352.66"
$ws.Range("E21").Value = "  -7.45%  "
$ws.Range("D22").Value = "'6.52"
$ws.Range("E22").Value = "  -3.00%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'5.72"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").Value = "'64.97"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").Value = "'0.176"
$ws.Range("E27").Value = "  -6.85%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -3.81%  "
$ws.Range("D30").Value = "0.0₃0854"
$ws.Range("E30").Value = "  -9.39%  "
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "'19.62"
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("D34").Value = "'152.58"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").Value = "'5.56"
$ws.Range("E36").Value = "  -6.15%  "
$ws.Range("E37").Value = "  -7.20%  "
$ws.Range("E38").Value = "  -5.90%  "
$ws.Range("D39").Value = "'37.52"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("D41").Value = "'3.71"
$ws.Range("E41").Value = "  -5.06%  "
$ws.Range("D42").Value = "2.287.52"
$ws.Range("E42").Value = "  -5.47%  "
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "'20.32"
$ws.Range("E45").Value = "  -7.85%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "'4.94"
$ws.Range("E47").Value = "  -4.46%  "
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").Value = "'18.40"
$ws.Range("E51").Value = "  -7.27%  "
